# Auto-generated edit script: apply numeric corrections to H:N columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR job-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 624.1667
$ws.Range("I8").Value = 33.714287
$ws.Range("J8").Value = 1450.8
$ws.Range("K8").Value = 101.142861
$ws.Range("L8").Value = 4352.4
$ws.Range("M8").Value = 37.857139
$ws.Range("N8").Value = -4630.4

$ws.Range("H9").Value = 152.71428
$ws.Range("I9").Value = 133.8
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 133.8
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 35.19999999999999
$ws.Range("N9").Value = -538

$ws.Range("H28").Value = 823.85187
$ws.Range("I28").Value = 914.8570999999999
$ws.Range("J28").Value = 725.8461
$ws.Range("K28").Value = 914.8570999999999
$ws.Range("L28").Value = 725.8461
$ws.Range("M28").Value = -429.8570999999999
$ws.Range("N28").Value = -1695.8461

$ws.Range("H62").Value = 4800
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -6048

$ws.Range("H65").Value = 4800
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -30240

$ws.Range("H96").Value = 13148.25
$ws.Range("I96").Value = 725.8
$ws.Range("J96").Value = 33852.332
$ws.Range("K96").Value = 2177.4
$ws.Range("L96").Value = 101556.996
$ws.Range("M96").Value = -804.3999999999996
$ws.Range("N96").Value = -104302.996

$ws.Range("H98").Value = 3066.889
$ws.Range("I98").Value = 2680.2666
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2680.2666
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -1182.2666
$ws.Range("N98").Value = -7996

$ws.Range("H113").Value = 2272.7273
$ws.Range("I113").Value = 2666.6667
$ws.Range("J113").Value = 2125
$ws.Range("K113").Value = 2666.6667
$ws.Range("L113").Value = 2125
$ws.Range("M113").Value = 587.3332999999998
$ws.Range("N113").Value = -8633

$ws.Range("H122").Value = 3066.889
$ws.Range("I122").Value = 2680.2666
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8040.7998
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5590.7998
$ws.Range("N122").Value = -19900

$ws.Range("H127").Value = 1608.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1608.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 4825.5
$ws.Range("M127").Value = ""
$ws.Range("N127").Value = -14745.5

$ws.Range("H129").Value = 1434.8064
$ws.Range("I129").Value = 763.9167
$ws.Range("J129").Value = 1858.5264
$ws.Range("K129").Value = 2291.7501
$ws.Range("L129").Value = 5575.5792
$ws.Range("M129").Value = 2708.2499
$ws.Range("N129").Value = -15575.5792

$ws.Range("H137").Value = 43488344
$ws.Range("I137").Value = 992.5
$ws.Range("K137").Value = 2977.5
$ws.Range("M137").Value = -427.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 334837
$ws.Range("I45").Value = 501005.5
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 501005.5
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -500628.5
$ws.Range("N45").Value = -3254

$ws.Range("H102").Value = 1803.0588
$ws.Range("I102").Value = 1702.5
$ws.Range("J102").Value = 2044.4
$ws.Range("K102").Value = 1702.5
$ws.Range("L102").Value = 2044.4
$ws.Range("M102").Value = -80.5
$ws.Range("N102").Value = -5288.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1627.2727
$ws.Range("I105").Value = 1661.125
$ws.Range("J105").Value = 1537
$ws.Range("K105").Value = 1661.125
$ws.Range("L105").Value = 1537
$ws.Range("M105").Value = 85.875
$ws.Range("N105").Value = -5031

$ws.Range("H134").Value = 53595.332
$ws.Range("I134").Value = 53595.332
$ws.Range("K134").Value = 160785.996
$ws.Range("M134").Value = -158250.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1492.2
$ws.Range("I16").Value = 1435.7778
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1435.7778
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1148.7778
$ws.Range("N16").Value = -2574

$ws.Range("H31").Value = 2746.4092
$ws.Range("I31").Value = 1361.4
$ws.Range("K31").Value = 1361.4
$ws.Range("M31").Value = -1066.4

$ws.Range("H34").Value = 2746.4092
$ws.Range("I34").Value = 1361.4
$ws.Range("K34").Value = 1361.4
$ws.Range("M34").Value = -1159.4

$ws.Range("H99").Value = 1311.3334
$ws.Range("I99").Value = 1310
$ws.Range("J99").Value = 1314
$ws.Range("K99").Value = 1310
$ws.Range("L99").Value = 1314
$ws.Range("M99").Value = 188
$ws.Range("N99").Value = -4310

$ws.Range("H113").Value = 1492.2
$ws.Range("I113").Value = 1435.7778
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1435.7778
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 734.2221999999999
$ws.Range("N113").Value = -6340

$ws.Range("H126").Value = 1311.3334
$ws.Range("I126").Value = 1310
$ws.Range("J126").Value = 1314
$ws.Range("K126").Value = 3930
$ws.Range("L126").Value = 3942
$ws.Range("M126").Value = -1460
$ws.Range("N126").Value = -8882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 70707190
$ws.Range("J96").Value = 70707190
$ws.Range("L96").Value = 212121570
$ws.Range("N96").Value = -212125688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2866.5
$ws.Range("I102").Value = 2488.6667
$ws.Range("K102").Value = 2488.6667
$ws.Range("M102").Value = -866.6667000000002

$ws.Range("H107").Value = 531.13635
$ws.Range("I107").Value = 424.23077
$ws.Range("J107").Value = 685.55554
$ws.Range("K107").Value = 424.23077
$ws.Range("L107").Value = 685.55554
$ws.Range("M107").Value = 1495.76923
$ws.Range("N107").Value = -4525.55554

$ws.Range("H113").Value = 1633.6666
$ws.Range("I113").Value = 1575.375
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1575.375
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 594.625
$ws.Range("N113").Value = -6440

$ws.Range("H126").Value = 2660.1667
$ws.Range("I126").Value = 2902.2
$ws.Range("J126").Value = 1450
$ws.Range("K126").Value = 8706.599999999999
$ws.Range("L126").Value = 4350
$ws.Range("M126").Value = -6236.599999999999
$ws.Range("N126").Value = -9290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1525.9565
$ws.Range("I7").Value = 1385.5714
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 1385.5714
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -1273.5714
$ws.Range("N7").Value = -3224

$ws.Range("H61").Value = 2581
$ws.Range("I61").Value = 2966.6667
$ws.Range("J61").Value = 2002.5
$ws.Range("K61").Value = 2966.6667
$ws.Range("L61").Value = 2002.5
$ws.Range("M61").Value = -2764.6667
$ws.Range("N61").Value = -2406.5

$ws.Range("H100").Value = 3033.111
$ws.Range("I100").Value = 1400
$ws.Range("J100").Value = 5074.5
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 5074.5
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -6156.5

$ws.Range("H104").Value = 17000
$ws.Range("J104").Value = 17000
$ws.Range("L104").Value = 17000
$ws.Range("N104").Value = -23988

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

$ws.Range("H106").Value = 21974.75
$ws.Range("J106").Value = 21974.75
$ws.Range("L106").Value = 21974.75
$ws.Range("N106").Value = -24498.75

$ws.Range("H113").Value = 2581
$ws.Range("I113").Value = 2966.6667
$ws.Range("J113").Value = 2002.5
$ws.Range("K113").Value = 2966.6667
$ws.Range("L113").Value = 2002.5
$ws.Range("M113").Value = -796.6667000000002
$ws.Range("N113").Value = -6342.5

$ws.Range("H126").Value = 1525.9565
$ws.Range("I126").Value = 1385.5714
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4156.7142
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1686.7142
$ws.Range("N126").Value = -13940

$ws.Range("H136").Value = 1986.85
$ws.Range("I136").Value = 1677.2858
$ws.Range("J136").Value = 2153.5386
$ws.Range("K136").Value = 5031.857400000001
$ws.Range("L136").Value = 6460.6158
$ws.Range("M136").Value = -2481.857400000001
$ws.Range("N136").Value = -11560.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 23333.334
$ws.Range("I19").Value = 30000
$ws.Range("J19").Value = 20000
$ws.Range("K19").Value = 30000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = -29826
$ws.Range("N19").Value = -20348

$ws.Range("H113").Value = 477.9643
$ws.Range("I113").Value = 374.4
$ws.Range("J113").Value = 597.46155
$ws.Range("K113").Value = 1123.2
$ws.Range("L113").Value = 1792.38465
$ws.Range("M113").Value = 1046.8
$ws.Range("N113").Value = -6132.38465

$ws.Range("H126").Value = 1832.96
$ws.Range("I126").Value = 1764.7273
$ws.Range("K126").Value = 5294.1819
$ws.Range("M126").Value = -2824.1819

$ws.Range("H136").Value = 780.5641000000001
$ws.Range("I136").Value = 679.125
$ws.Range("J136").Value = 942.86664
$ws.Range("K136").Value = 2037.375
$ws.Range("L136").Value = 2828.59992
$ws.Range("M136").Value = 512.625
$ws.Range("N136").Value = -7928.59992
